$d = $word.ActiveDocument

# The run's text changes from "1999 to 2000" to "2000", and the run is
# moved so it now comes *after* the "_GoBack" bookmark instead of before it.
# Do the text edit first, then relocate the bookmark around the run.

$d.Content.Find.Execute("1999 to 2000", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2000", 2)

# Remove the existing "_GoBack" bookmark (currently sits right after the run).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Locate the (now shortened) run's text and re-create the bookmark as a
# zero-length range immediately in front of it, so bookmarkStart/bookmarkEnd
# precede the run in document order.
$rng = $d.Content
$rng.Find.Execute("2000", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$start = $rng.Start
$bmRange = $d.Range($start, $start)
$d.Bookmarks.Add("_GoBack", $bmRange)
